$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H4").Value = 191.5
$ws_ALC.Range("I4").Value = 137.75
$ws_ALC.Range("K4").Value = 137.75
$ws_ALC.Range("M4").Value = -23.75

$ws_ALC.Range("H28").Value = 1363.1875
$ws_ALC.Range("J28").Value = 1381.6666
$ws_ALC.Range("L28").Value = 1381.6666
$ws_ALC.Range("N28").Value = -2351.6666

$ws_ALC.Range("H40").Value = 5411.1113
$ws_ALC.Range("I40").Value = 4322.222
$ws_ALC.Range("J40").Value = 6500
$ws_ALC.Range("K40").Value = 4322.222
$ws_ALC.Range("L40").Value = 6500
$ws_ALC.Range("M40").Value = -4147.222
$ws_ALC.Range("N40").Value = -6850

$ws_ALC.Range("H42").Value = 1447.2
$ws_ALC.Range("I42").Value = 177.875
$ws_ALC.Range("J42").Value = 6524.5
$ws_ALC.Range("K42").Value = 533.625
$ws_ALC.Range("L42").Value = 19573.5
$ws_ALC.Range("M42").Value = -303.625
$ws_ALC.Range("N42").Value = -20033.5

$ws_ALC.Range("H51").Value = 7642.6665
$ws_ALC.Range("J51").Value = 7747.5
$ws_ALC.Range("L51").Value = 7747.5
$ws_ALC.Range("N51").Value = -8715.5

$ws_ALC.Range("H70").Value = 3650
$ws_ALC.Range("J70").Value = 3650
$ws_ALC.Range("L70").Value = 10950
$ws_ALC.Range("N70").Value = -11490

$ws_ALC.Range("H73").Value = 3650
$ws_ALC.Range("J73").Value = 3650
$ws_ALC.Range("L73").Value = 10950
$ws_ALC.Range("N73").Value = -12822

$ws_ALC.Range("H80").Value = 2579.4075
$ws_ALC.Range("I80").Value = 607.75
$ws_ALC.Range("J80").Value = 4156.7334
$ws_ALC.Range("K80").Value = 1823.25
$ws_ALC.Range("L80").Value = 12470.2002
$ws_ALC.Range("M80").Value = -825.25
$ws_ALC.Range("N80").Value = -14466.2002

$ws_ALC.Range("H83").Value = 2579.4075
$ws_ALC.Range("I83").Value = 607.75
$ws_ALC.Range("J83").Value = 4156.7334
$ws_ALC.Range("K83").Value = 5469.75
$ws_ALC.Range("L83").Value = 37410.6006
$ws_ALC.Range("M83").Value = -477.75
$ws_ALC.Range("N83").Value = -47394.6006

$ws_ALC.Range("H98").Value = 1057.9286
$ws_ALC.Range("I98").Value = 754.2308
$ws_ALC.Range("K98").Value = 754.2308
$ws_ALC.Range("M98").Value = 743.7692

$ws_ALC.Range("H107").Value = 1215
$ws_ALC.Range("I107").Value = 1215
$ws_ALC.Range("K107").Value = 1215
$ws_ALC.Range("M107").Value = 705

$ws_ALC.Range("H122").Value = 1057.9286
$ws_ALC.Range("I122").Value = 754.2308
$ws_ALC.Range("K122").Value = 2262.6924
$ws_ALC.Range("M122").Value = 187.3076000000001

$ws_ALC.Range("H125").Value = 1690.75
$ws_ALC.Range("I125").Value = 1538
$ws_ALC.Range("K125").Value = 13842
$ws_ALC.Range("M125").Value = -11382

$ws_ALC.Range("H132").Value = 17856.906
$ws_ALC.Range("I132").Value = 1757.8928
$ws_ALC.Range("K132").Value = 5273.678400000001
$ws_ALC.Range("M132").Value = -2743.678400000001

$ws_ALC.Range("H135").Value = 20333.334
$ws_ALC.Range("I135").Value = 0
$ws_ALC.Range("K135").Value = 0
$ws_ALC.Range("M135").ClearContents()

$ws_ALC.Range("H141").Value = 8868.857
$ws_ALC.Range("I141").Value = 8529
$ws_ALC.Range("K141").Value = 25587
$ws_ALC.Range("M141").Value = -20407

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 8344.263000000001
$ws_ARM.Range("I61").Value = 1616.75
$ws_ARM.Range("K61").Value = 1616.75
$ws_ARM.Range("M61").Value = -1404.75

$ws_ARM.Range("H74").Value = 5028.3076
$ws_ARM.Range("I74").Value = 4926.8
$ws_ARM.Range("K74").Value = 4926.8
$ws_ARM.Range("M74").Value = -4052.8

$ws_ARM.Range("H77").Value = 5028.3076
$ws_ARM.Range("I77").Value = 4926.8
$ws_ARM.Range("K77").Value = 24634
$ws_ARM.Range("M77").Value = -20266

$ws_ARM.Range("H122").Value = 4997.3076
$ws_ARM.Range("I122").Value = 4370.625
$ws_ARM.Range("K122").Value = 13111.875
$ws_ARM.Range("M122").Value = -10661.875

$ws_ARM.Range("H132").Value = 2036.037
$ws_ARM.Range("I132").Value = 1494.8235
$ws_ARM.Range("K132").Value = 4484.470499999999
$ws_ARM.Range("M132").Value = -1954.470499999999

$ws_ARM.Range("H136").Value = 8344.263000000001
$ws_ARM.Range("I136").Value = 1616.75
$ws_ARM.Range("K136").Value = 4850.25
$ws_ARM.Range("M136").Value = -2300.25

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 4724.727
$ws_BSM.Range("I94").Value = 1118.75
$ws_BSM.Range("J94").Value = 6785.2856
$ws_BSM.Range("K94").Value = 1118.75
$ws_BSM.Range("L94").Value = 6785.2856
$ws_BSM.Range("M94").Value = -667.75
$ws_BSM.Range("N94").Value = -7687.2856

$ws_BSM.Range("H100").Value = 15178.4
$ws_BSM.Range("J100").Value = 15178.4
$ws_BSM.Range("L100").Value = 15178.4
$ws_BSM.Range("N100").Value = -17342.4

$ws_BSM.Range("H105").Value = 3408
$ws_BSM.Range("I105").Value = 2174.5
$ws_BSM.Range("K105").Value = 2174.5
$ws_BSM.Range("M105").Value = -427.5

$ws_BSM.Range("H107").Value = 10814.571
$ws_BSM.Range("I107").Value = 9685
$ws_BSM.Range("K107").Value = 9685
$ws_BSM.Range("M107").Value = -7765

$ws_BSM.Range("H134").Value = 2065.9285
$ws_BSM.Range("I134").Value = 1539.091
$ws_BSM.Range("K134").Value = 4617.272999999999
$ws_BSM.Range("M134").Value = -2082.272999999999

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H58").Value = 2948.5789
$ws_CRP.Range("I58").Value = 1965.6666
$ws_CRP.Range("J58").Value = 3833.2
$ws_CRP.Range("K58").Value = 1965.6666
$ws_CRP.Range("L58").Value = 3833.2
$ws_CRP.Range("M58").Value = -1762.6666
$ws_CRP.Range("N58").Value = -4239.2

$ws_CRP.Range("H99").Value = 1500
$ws_CRP.Range("I99").Value = 1500
$ws_CRP.Range("J99").Value = 0
$ws_CRP.Range("K99").Value = 1500
$ws_CRP.Range("L99").Value = 0
$ws_CRP.Range("M99").Value = -2
$ws_CRP.Range("N99").ClearContents()

$ws_CRP.Range("H122").Value = 4539.1333
$ws_CRP.Range("I122").Value = 3933.7
$ws_CRP.Range("J122").Value = 5750
$ws_CRP.Range("K122").Value = 11801.1
$ws_CRP.Range("L122").Value = 17250
$ws_CRP.Range("M122").Value = -9351.099999999999
$ws_CRP.Range("N122").Value = -22150

$ws_CRP.Range("H126").Value = 1500
$ws_CRP.Range("I126").Value = 1500
$ws_CRP.Range("J126").Value = 0
$ws_CRP.Range("K126").Value = 4500
$ws_CRP.Range("L126").Value = 0
$ws_CRP.Range("M126").Value = -2030
$ws_CRP.Range("N126").ClearContents()

$ws_CRP.Range("H132").Value = 2833.1177
$ws_CRP.Range("I132").Value = 2871.8125
$ws_CRP.Range("K132").Value = 8615.4375
$ws_CRP.Range("M132").Value = -6085.4375

$ws_CRP.Range("H136").Value = 2948.5789
$ws_CRP.Range("I136").Value = 1965.6666
$ws_CRP.Range("J136").Value = 3833.2
$ws_CRP.Range("K136").Value = 5896.9998
$ws_CRP.Range("L136").Value = 11499.6
$ws_CRP.Range("M136").Value = -3346.9998
$ws_CRP.Range("N136").Value = -16599.6

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 26.105263
$ws_CUL.Range("J2").Value = 22.846153
$ws_CUL.Range("L2").Value = 137.076918
$ws_CUL.Range("N2").Value = -363.076918

$ws_CUL.Range("H38").Value = 679.4595
$ws_CUL.Range("J38").Value = 1465.875
$ws_CUL.Range("L38").Value = 4397.625
$ws_CUL.Range("N38").Value = -5091.625

$ws_CUL.Range("H81").Value = 0
$ws_CUL.Range("I81").Value = 0
$ws_CUL.Range("J81").Value = 0
$ws_CUL.Range("K81").Value = 0
$ws_CUL.Range("L81").Value = 0
$ws_CUL.Range("M81").ClearContents()
$ws_CUL.Range("N81").ClearContents()

$ws_CUL.Range("H84").Value = 0
$ws_CUL.Range("I84").Value = 0
$ws_CUL.Range("J84").Value = 0
$ws_CUL.Range("K84").Value = 0
$ws_CUL.Range("L84").Value = 0
$ws_CUL.Range("M84").ClearContents()
$ws_CUL.Range("N84").ClearContents()

$ws_CUL.Range("H131").Value = 1669695.9
$ws_CUL.Range("I131").Value = 4000633.5
$ws_CUL.Range("K131").Value = 12001900.5
$ws_CUL.Range("M131").Value = -11996860.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 2695.2727
$ws_GSM.Range("I102").Value = 2695.2727
$ws_GSM.Range("K102").Value = 2695.2727
$ws_GSM.Range("M102").Value = -1073.2727

$ws_GSM.Range("H122").Value = 4953.857
$ws_GSM.Range("I122").Value = 4355.467
$ws_GSM.Range("K122").Value = 13066.401
$ws_GSM.Range("M122").Value = -10616.401

$ws_GSM.Range("H132").Value = 6882.614
$ws_GSM.Range("I132").Value = 6411.273
$ws_GSM.Range("K132").Value = 19233.819
$ws_GSM.Range("M132").Value = -16703.819

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H16").Value = 12500399
$ws_LTW.Range("J16").Value = 1150
$ws_LTW.Range("L16").Value = 1150
$ws_LTW.Range("N16").Value = -1490

$ws_LTW.Range("H40").Value = 200000
$ws_LTW.Range("I40").Value = 0
$ws_LTW.Range("J40").Value = 200000
$ws_LTW.Range("K40").Value = 0
$ws_LTW.Range("L40").Value = 200000
$ws_LTW.Range("M40").ClearContents()
$ws_LTW.Range("N40").Value = -200272

$ws_LTW.Range("H46").Value = 2844.8462
$ws_LTW.Range("I46").Value = 1159.8
$ws_LTW.Range("K46").Value = 1159.8
$ws_LTW.Range("M46").Value = -971.8

$ws_LTW.Range("H55").Value = 571.2778
$ws_LTW.Range("J55").Value = 712.5
$ws_LTW.Range("L55").Value = 712.5
$ws_LTW.Range("N55").Value = -1058.5

$ws_LTW.Range("H123").Value = 39750
$ws_LTW.Range("J123").Value = 39750
$ws_LTW.Range("L123").Value = 39750
$ws_LTW.Range("N123").Value = -49550

$ws_LTW.Range("H136").Value = 2321.8076
$ws_LTW.Range("I136").Value = 2358.4443
$ws_LTW.Range("J136").Value = 2239.375
$ws_LTW.Range("K136").Value = 7075.3329
$ws_LTW.Range("L136").Value = 6718.125
$ws_LTW.Range("M136").Value = -4525.3329
$ws_LTW.Range("N136").Value = -11818.125

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H57").Value = 29909.092
$ws_WVR.Range("J57").Value = 29909.092
$ws_WVR.Range("L57").Value = 29909.092
$ws_WVR.Range("N57").Value = -31417.092

$ws_WVR.Range("H97").Value = 40000
$ws_WVR.Range("J97").Value = 40000
$ws_WVR.Range("L97").Value = 40000
$ws_WVR.Range("N97").Value = -41982

$ws_WVR.Range("H122").Value = 1713.9
$ws_WVR.Range("I122").Value = 1713.9
$ws_WVR.Range("K122").Value = 5141.700000000001
$ws_WVR.Range("M122").Value = -2691.700000000001
